# Update CVD (ytd) figures + a handful of monthly cells across several
# location sheets, and insert a missing "Internal Fill Rate" /
# "Commit/Forecast" row on the Lincoln Missouri sheet so the monthly CVD
# data stays dynamic (per commit message).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Tipp City Ohio: Professional Voluntary Turnover / Commit-Forecast, Apr
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tipp City Ohio")
$ws.Range("L4").Value = 0.3636

# ---------------------------------------------------------------------
# Milwaukee Pmc Hq Wisconsin: Internal Fill Rate / Commit-Forecast, Apr
# clear the cell entirely (was an explicit 0)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("L3").ClearContents()

# ---------------------------------------------------------------------
# Piedras Negras Fasco Mexico: ytd (cvd) column updated for both rows
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Piedras Negras Fasco Mexico")
$ws.Range("E2").Value = 0.0776
$ws.Range("E3").Value = 0.0776

# ---------------------------------------------------------------------
# Faridabad India: ytd (cvd) column updated for rows 5-7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Faridabad India")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# ---------------------------------------------------------------------
# Fort Wayne Indiana: Internal Fill Rate / Commit-Forecast monthly figures
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fort Wayne Indiana")
$ws.Range("G4").Value = 0.0263
$ws.Range("I4").Value = 0.027
$ws.Range("J4").Value = 0.0531
$ws.Range("L4").Value = 0.0278

# ---------------------------------------------------------------------
# Grafton Wisconsin: ytd (cvd) column updated for rows 5-6
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Grafton Wisconsin")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776

# ---------------------------------------------------------------------
# Hyderabad India: Professional Voluntary Turnover / Commit-Forecast, Apr
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")
$ws.Range("L4").Value = 0.0172

# ---------------------------------------------------------------------
# Lincoln Missouri: ytd for the two existing "Internal Fill Rate" rows,
# then insert the missing "Internal Fill Rate" / "Commit/Forecast" row
# (pushing the Manufacturing Voluntary Turnover block down one row), and
# refresh the ytd + monthly figures on the shifted rows.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lincoln Missouri")
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "PES"
$ws.Range("B7").Value = "PES NA Motors Solutions"
$ws.Range("C7").Value = "Lincoln Missouri"
$ws.Range("D7").Value = "Internal Fill Rate"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "Commit/Forecast"
$ws.Range("L7:N7").Value = 0
$ws.Range("O7:W7").Value = 0

# row 8 == old row 7 (Manufacturing Voluntary Turnover / PY Actual)
$ws.Range("E8").Value = 0.0776

# row 9 == old row 8 (Manufacturing Voluntary Turnover / AOP)
$ws.Range("E9").Value = 0.0776

# row 10 == old row 9 (Manufacturing Voluntary Turnover / Commit-Forecast)
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0.0349

# ---------------------------------------------------------------------
# Piedras Negras Jakel Mexico: clear L5, update ytd + monthly figures
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Piedras Negras Jakel Mexico")
$ws.Range("L5").ClearContents()
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("G8").Value = 0.0289
$ws.Range("J8").Value = 0.0522
$ws.Range("K8").Value = 0.0321
$ws.Range("L8").Value = 0.021
